$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns to reflect the latest scrape.
# Numeric-looking Price values are written with a leading quote-prefix so Excel
# keeps them as text (matching the original "inlineStr" formatting, e.g. trailing
# zeros like "1.00" or dotted thousands like "1.895.92"), then the style is reset
# back to Normal so no stray number-format styling is introduced.

$ws.Range('D2').Value = '27.542.24'
$ws.Range('E2').Value = '  -3.24%  '
$ws.Range('D3').Value = '1.660.61'
$ws.Range('E3').Value = '  -3.72%  '
$ws.Range('D4').Value = "'1.00"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = "'214.54"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.96%  '
$ws.Range('E6').Value = '  -2.26%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').Value = "'24.34"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.72%  '
$ws.Range('D9').Value = "'0.263"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.66%  '
$ws.Range('E10').Value = '  -2.33%  '
$ws.Range('E11').Value = '  -1.59%  '
$ws.Range('D12').Value = '1.895.92'
$ws.Range('E12').Value = '  -3.69%  '
$ws.Range('D13').Value = '1.656.23'
$ws.Range('E13').Value = '  -4.01%  '
$ws.Range('E14').Value = '  -2.32%  '
$ws.Range('D15').Value = "'0.567"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.69%  '
$ws.Range('D16').Value = "'65.84"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.60%  '
$ws.Range('D17').Value = '27.554.72'
$ws.Range('E17').Value = '  -2.90%  '
$ws.Range('D18').Value = "'239.48"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.15%  '
$ws.Range('E19').Value = '  -2.66%  '
$ws.Range('D20').Value = "'7.68"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.34%  '
$ws.Range('E21').Value = '  +0.05%  '
$ws.Range('E22').Value = '  -3.09%  '
$ws.Range('D23').Value = "'9.41"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E24').Value = '  -1.21%  '
$ws.Range('D25').Value = "'146.10"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.02%  '
$ws.Range('D26').Value = "'7.22"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.73%  '
$ws.Range('D27').Value = "'16.23"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.29%  '
$ws.Range('D28').Value = "'1.00"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.17%  '
$ws.Range('D29').Value = "'0.111"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.20%  '
$ws.Range('D30').Value = "'0.0501"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.99%  '
$ws.Range('E31').Value = '  -0.62%  '
$ws.Range('E32').Value = '  -2.86%  '
$ws.Range('D33').Value = '1.453.50'
$ws.Range('E33').Value = '  -1.94%  '
$ws.Range('E34').Value = '  -4.04%  '
$ws.Range('E35').Value = '  -4.10%  '
$ws.Range('D36').Value = "'2.39"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.74%  '
$ws.Range('D37').Value = "'0.920"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -5.69%  '
$ws.Range('D38').Value = "'0.573"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.40%  '
$ws.Range('E39').Value = '  -2.72%  '
$ws.Range('E40').Value = '  +0.29%  '
$ws.Range('E41').Value = '  -0.04%  '
$ws.Range('D42').Value = "'66.56"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.37%  '
$ws.Range('E43').Value = '  -2.88%  '
$ws.Range('E44').Value = '  -2.69%  '
$ws.Range('E45').Value = '  -3.66%  '
$ws.Range('E46').Value = '  -2.03%  '
$ws.Range('E47').Value = '  -0.13%  '
$ws.Range('D48').Value = "'88.77"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.81%  '
$ws.Range('D49').Value = '0.0₆0107'
$ws.Range('E49').Value = '  -5.13%  '
$ws.Range('E50').Value = '  -1.29%  '
$ws.Range('D51').Value = "'7.85"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.93%  '
